# Auto-generated edit script: update cryptos list values per upstream diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '89.557.99'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  -1.56%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.081.87'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  -2.32%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '235.18'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +8.61%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '617.92'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -1.30%  '
$ws.Range('E7').Value = '  -6.31%  '
$ws.Range('E8').Value = '  -1.63%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.00'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +0.05%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '3.080.49'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -2.33%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.713'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -6.11%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.198'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -1.57%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000251'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +1.93%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '35.21'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +0.32%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '89.350.98'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -1.43%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.36'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -6.47%  '
$ws.Range('E17').Value = '  -2.61%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.089.84'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -3.38%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.78'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +0.08%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0000211'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -0.20%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.76'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -5.70%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '431.87'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -9.12%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.38'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +2.11%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '8.74'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -4.60%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '5.58'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -6.14%  '
$ws.Range('B26').Value = 'Litecoin'
$ws.Range('C26').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '87.01'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -8.29%  '
$ws.Range('B27').Value = 'Aptos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.70'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -5.16%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '3.250.07'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -2.21%  '
$ws.Range('E29').Value = '  +0.40%  '
$ws.Range('E30').Value = '  +14.76%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '9.03'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -3.10%  '
$ws.Range('E32').Value = '  -4.52%  '
$ws.Range('E33').Value = '  -10.98%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '25.56'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -6.32%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.149'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +2.05%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '7.12'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +1.29%  '
$ws.Range('B37').Value = 'Bittensor'
$ws.Range('C37').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '494.91'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -5.04%  '
$ws.Range('B38').Value = 'dogwifhat'
$ws.Range('C38').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.65'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +0.37%  '
$ws.Range('E39').Value = '  -3.44%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.25'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -5.29%  '
$ws.Range('B41').Value = 'MantraDAO'
$ws.Range('C41').Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.63'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +55.25%  '
$ws.Range('B42').Value = 'Hedera'
$ws.Range('C42').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0891'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -1.30%  '
$ws.Range('E43').Value = '  -0.63%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.397'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -7.50%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '152.75'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +1.57%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.84'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -7.17%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.674'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -8.25%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '44.47'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -1.97%  '
$ws.Range('B50').Value = 'FirstDigitalUSD'
$ws.Range('C50').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.00'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -0.08%  '
$ws.Range('B51').Value = 'ImmutableX'
$ws.Range('C51').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.30'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -5.20%  '
